# Natmi following Dr Hou advice
# Update the LR-pair stats for Ifnk-Ifnar1 (rows 2-4) to reflect the new
# ligand/receptor expressing-cell counts (1 -> 3) and the resulting
# recalculated expression / specificity values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2 = @{
        E = 3
        G = 0.8715449999999999
        H = 2.614635
        K = 3
        M = 24.577204
        N = 73.731612
        O = 0.5324151489760768
        P = 0.5324151489760768
        Q = 21.42013926018
        R = 192.78125334162
        S = 0.5324151489760768
        T = 0.5324151489760768
    }
    3 = @{
        E = 3
        G = 0.8715449999999999
        H = 2.614635
        K = 3
        M = 15.033452
        N = 45.100356
        O = 0.3256691683156758
        P = 0.3256691683156758
        Q = 13.10232992334
        R = 117.92096931006
        S = 0.3256691683156758
        T = 0.3256691683156758
    }
    4 = @{
        E = 3
        G = 0.8715449999999999
        H = 2.614635
        K = 3
        M = 6.551073333333334
        N = 19.65322
        O = 0.1419156827082475
        P = 0.1419156827082475
        Q = 5.7095552083
        R = 51.3859968747
        S = 0.1419156827082475
        T = 0.1419156827082475
    }
}

foreach ($row in $updates.Keys) {
    $cols = $updates[$row]
    foreach ($col in $cols.Keys) {
        $ws.Range("$col$row").Value = $cols[$col]
    }
}
